$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - OP1_Wide -AES
$ws.Range("C2").Value = 415000000.0000029
$ws.Range("D2").Value = 129.5714928739288
$ws.Range("E2").Value = 53772169542.68082

# Row 3 - OP1_Wide -ISAGEN
$ws.Range("C3").Value = 459458287.0000031
$ws.Range("D3").Value = 81.67845026912242
$ws.Range("E3").Value = 37527840845.46593

# Row 4 - OP1_Wide- EPM
$ws.Range("C4").Value = 880608868
$ws.Range("D4").Value = 53.76483966340702
$ws.Range("E4").Value = 47345794594.19435

# Row 5 - TOTAL / TODAS LAS OFERTAS
$ws.Range("C5").Value = 1755067155.000006
$ws.Range("D5").Value = 78.9974358458897
$ws.Range("E5").Value = 138645804982.3411
